$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) First paragraph: append "  (This is a change – Version for main
#    branch)" after the existing text, with the appended text in red
#    and split across three runs (matching how Word recorded the
#    original edit).
# ------------------------------------------------------------------
$firstPara = $d.Paragraphs.Item(1)
$paraStart = $firstPara.Range.Start
$paraTextEnd = $firstPara.Range.End - 1   # exclude the paragraph mark

$baseRange = $d.Range($paraStart, $paraTextEnd)
$baseRange.InsertAfter("  ")
$pos = $paraTextEnd + 2

$part1 = "(This is a change " + [char]0x2013 + " Ve"
$part2 = "rsion for main branch"
$part3 = ")"

$ins1 = $d.Range($pos, $pos)
$ins1.InsertAfter($part1)
$fmt1 = $d.Range($pos, $pos + $part1.Length)
$fmt1.Font.Color = 255
$pos = $pos + $part1.Length

$ins2 = $d.Range($pos, $pos)
$ins2.InsertAfter($part2)
$fmt2 = $d.Range($pos, $pos + $part2.Length)
$fmt2.Font.Color = 255
$pos = $pos + $part2.Length

$ins3 = $d.Range($pos, $pos)
$ins3.InsertAfter($part3)
$fmt3 = $d.Range($pos, $pos + $part3.Length)
$fmt3.Font.Color = 255
$pos = $pos + $part3.Length

# ------------------------------------------------------------------
# 2) Remove the final paragraph ("ank God almighty, we are free at
#    last.") entirely.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.Delete()
